$wb = $excel.ActiveWorkbook

# --- 1. Update localization status text: "Ready for handoff" -> "In Translation" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value2 = "In Translation"
$overview.Range("F2").Value2 = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value2 = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value2 = "In Translation"

# --- 2. Shrink the now-narrower "Status" columns to fit the shorter text ---
# (Ready for handoff -> In Translation made these columns narrower in the
# regenerated report.) ColumnWidth is expressed in characters; 12.5 is the
# closest the engine's character-grid will resolve to the target width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
